$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 552.125
$ws.Range("J92").Value = 785
$ws.Range("L92").Value = 785
$ws.Range("N92").Value = -3281
$ws.Range("H100").Value = 5907.615
$ws.Range("I100").Value = 5537.6665
$ws.Range("K100").Value = 5537.6665
$ws.Range("M100").Value = -4996.6665
$ws.Range("H112").Value = 1221.97
$ws.Range("J112").Value = 1295.6195
$ws.Range("L112").Value = 3886.8585
$ws.Range("N112").Value = -6102.8585
$ws.Range("H137").Value = 4352217.5
$ws.Range("I137").Value = 5267342
$ws.Range("J137").Value = 5375
$ws.Range("K137").Value = 15802026
$ws.Range("L137").Value = 16125
$ws.Range("M137").Value = -15799476
$ws.Range("N137").Value = -21225

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4412.9756
$ws.Range("I32").Value = 3252.972
$ws.Range("J32").Value = 11900.272
$ws.Range("K32").Value = 3252.972
$ws.Range("L32").Value = 11900.272
$ws.Range("M32").Value = -2965.972
$ws.Range("N32").Value = -12474.272
$ws.Range("H61").Value = 3215.7407
$ws.Range("I61").Value = 1564
$ws.Range("J61").Value = 3911.2104
$ws.Range("K61").Value = 1564
$ws.Range("L61").Value = 3911.2104
$ws.Range("M61").Value = -1352
$ws.Range("N61").Value = -4335.2104
$ws.Range("H74").Value = 740.4545000000001
$ws.Range("I74").Value = 740.4545000000001
$ws.Range("K74").Value = 740.4545000000001
$ws.Range("M74").Value = 133.5454999999999
$ws.Range("H77").Value = 740.4545000000001
$ws.Range("I77").Value = 740.4545000000001
$ws.Range("K77").Value = 3702.2725
$ws.Range("M77").Value = 665.7275
$ws.Range("H136").Value = 3215.7407
$ws.Range("I136").Value = 1564
$ws.Range("J136").Value = 3911.2104
$ws.Range("K136").Value = 4692
$ws.Range("L136").Value = 11733.6312
$ws.Range("M136").Value = -2142
$ws.Range("N136").Value = -16833.6312
$ws.Range("H137").Value = 44500
$ws.Range("J137").Value = 44500
$ws.Range("L137").Value = 44500
$ws.Range("N137").Value = -54700

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 39945
$ws.Range("J133").Value = 39945
$ws.Range("L133").Value = 39945
$ws.Range("N133").Value = -50065
$ws.Range("H134").Value = 2389.7812
$ws.Range("I134").Value = 1968.2174
$ws.Range("J134").Value = 3467.111
$ws.Range("K134").Value = 5904.6522
$ws.Range("L134").Value = 10401.333
$ws.Range("M134").Value = -3369.6522
$ws.Range("N134").Value = -15471.333
$ws.Range("H139").Value = 35000
$ws.Range("J139").Value = 35000
$ws.Range("L139").Value = 35000
$ws.Range("N139").Value = -45280

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2274726.8
$ws.Range("I31").Value = 2440816.5
$ws.Range("J31").Value = 4833.3335
$ws.Range("K31").Value = 2440816.5
$ws.Range("L31").Value = 4833.3335
$ws.Range("M31").Value = -2440521.5
$ws.Range("N31").Value = -5423.3335
$ws.Range("H34").Value = 2274726.8
$ws.Range("I34").Value = 2440816.5
$ws.Range("J34").Value = 4833.3335
$ws.Range("K34").Value = 2440816.5
$ws.Range("L34").Value = 4833.3335
$ws.Range("M34").Value = -2440614.5
$ws.Range("N34").Value = -5237.3335
$ws.Range("H74").Value = 23061.777
$ws.Range("J74").Value = 23061.777
$ws.Range("L74").Value = 23061.777
$ws.Range("N74").Value = -24809.777
$ws.Range("H77").Value = 23061.777
$ws.Range("J77").Value = 23061.777
$ws.Range("L77").Value = 69185.33099999999
$ws.Range("N77").Value = -77921.33099999999
$ws.Range("H132").Value = 2778.147
$ws.Range("I132").Value = 2187.318
$ws.Range("K132").Value = 6561.954000000001
$ws.Range("M132").Value = -4031.954000000001

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 6850.75
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 6850.75
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 20552.25
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -22424.25
$ws.Range("H83").Value = 6850.75
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 6850.75
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 61656.75
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -71016.75
$ws.Range("H122").Value = 1312.6875
$ws.Range("I122").Value = 651.5
$ws.Range("J122").Value = 1533.0834
$ws.Range("K122").Value = 5863.5
$ws.Range("L122").Value = 13797.7506
$ws.Range("M122").Value = -3413.5
$ws.Range("N122").Value = -18697.7506
$ws.Range("H132").Value = 2819.2727
$ws.Range("I132").Value = 1943.8572
$ws.Range("J132").Value = 4351.25
$ws.Range("K132").Value = 17494.7148
$ws.Range("L132").Value = 39161.25
$ws.Range("M132").Value = -14964.7148
$ws.Range("N132").Value = -44221.25

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 9800
$ws.Range("I96").Value = 9800
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 9800
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -7054
$ws.Range("N96").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1888937.5
$ws.Range("I136").Value = 3227675.8
$ws.Range("J136").Value = 2533.818
$ws.Range("K136").Value = 9683027.399999999
$ws.Range("L136").Value = 7601.454000000001
$ws.Range("M136").Value = -9680477.399999999
$ws.Range("N136").Value = -12701.454
$ws.Range("H140").Value = 29483.615
$ws.Range("J140").Value = 29483.615
$ws.Range("L140").Value = 29483.615
$ws.Range("N140").Value = -39843.61500000001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 22949.2
$ws.Range("I69").Value = 11246
$ws.Range("J69").Value = 25875
$ws.Range("K69").Value = 11246
$ws.Range("L69").Value = 25875
$ws.Range("M69").Value = -10497
$ws.Range("N69").Value = -27373
$ws.Range("H72").Value = 22949.2
$ws.Range("I72").Value = 11246
$ws.Range("J72").Value = 25875
$ws.Range("K72").Value = 33738
$ws.Range("L72").Value = 77625
$ws.Range("M72").Value = -29994
$ws.Range("N72").Value = -85113
$ws.Range("H95").Value = 28666.666
$ws.Range("J95").Value = 28666.666
$ws.Range("L95").Value = 28666.666
$ws.Range("N95").Value = -34158.666
$ws.Range("H122").Value = 591197.25
$ws.Range("I122").Value = 835679.4399999999
$ws.Range("K122").Value = 2507038.32
$ws.Range("M122").Value = -2504588.32
